$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 1.35
$ws.Range("B3").Value = 1.47
$ws.Range("F3").Value = 1.23
$ws.Range("E4").Value = 1.22
$ws.Range("D5").Value = 1.36
$ws.Range("F5").Value = 1.05
$ws.Range("C6").Value = 1.46
$ws.Range("E6").Value = 1.31
$ws.Range("F6").Value = 1.18
